$d = $word.ActiveDocument

# --- Step 1: split the "<Proposal Description>" paragraph, inserting a new
#     blank paragraph above it with the same paragraph/run formatting ---
$rng = $d.Content
$rng.Find.Execute("<Proposal Description>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$r = $para.Range
$r.Collapse(1)
$r.Text = "`r"

# --- Step 2: move the document's "_GoBack" bookmark from the signature
#     block onto the start of the (now second) "<Proposal Description>"
#     paragraph ---
$d.Bookmarks("_GoBack").Delete()

$rng2 = $d.Content
$rng2.Find.Execute("<Proposal Description>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $rng2.Paragraphs(1)
$r2 = $para2.Range
$r2.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r2)
